$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Pasquier Beignet chocolat" row (original row 384). Deleting the
# entire row shifts every subsequent row up by one, which is why the sheet's
# used range shrinks from A1:O392 to A1:O391.
$ws.Rows.Item(384).Delete()

# A handful of products toggled their "Online kein Bestand" (out-of-stock)
# label inside the aria-label text in column M. These rows are all above the
# deleted row, so their row numbers are unaffected by the deletion above.
$ws.Range("M12").Value = "Naturaplan Bio Vollkorntoast 10 Scheiben - Online kein Bestand 2.50 Schweizer Franken"
$ws.Range("M128").Value = "Pasquier Milchbrötchen 10St 3.80 Schweizer Franken"
$ws.Range("M243").Value = "Pasquier Schokobrötchen 16 Stück 8.95 Schweizer Franken"
$ws.Range("M275").Value = "Country Cracker Käse - Online kein Bestand 30% ab 2 Aktion 3.60 Schweizer Franken"
$ws.Range("M339").Value = "Pasquier Pancakes Choco 10 Stück 4.95 Schweizer Franken"

# The crawl was re-run later the same day, so every row's timestamp (column
# O) advances from 06:49:40 to 12:57:57. After the row deletion above, the
# data now spans rows 2 through 391.
for ($r = 2; $r -le 391; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-03-09 12:57:57"
}
